# Apply the two changes captured by the commit:
#   1. Slide 16's table switches from the deck's custom "Table_0" style
#      ({ADB0161C-4CCA-41CC-8C05-298B43189138}) to the built-in
#      "Medium Style 2 - Accent 1" style ({C5D83A1C-E755-4081-85AC-0FB3F9D9E480}).
#   2. The presentation's theme palette is swapped from the "Integral" theme
#      colours to the classic "Office Theme" colours (dk1/lt1 are identical
#      between the two palettes, so only the remaining 10 slots move).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{C5D83A1C-E755-4081-85AC-0FB3F9D9E480}")
    }
}

# --- 2. Recolour the theme from "Integral" to "Office Theme" ---------------
# dk1/lt1 (slots 1-2) are 000000/FFFFFF in both palettes, so they are left
# untouched; slots 3-12 move from the Integral values to the Office values.
$themeColors = $p.Slides.Item(1).ThemeColorScheme

$themeColors.Colors(3).RGB  = 6968388    # dk2      455F51 -> 44546A
$themeColors.Colors(4).RGB  = 15132391   # lt2      E3DED1 -> E7E6E6
$themeColors.Colors(5).RGB  = 13998939   # accent1  99CB38 -> 5B9BD5
$themeColors.Colors(6).RGB  = 3243501    # accent2  63A537 -> ED7D31
$themeColors.Colors(7).RGB  = 10855845   # accent3  E6D024 -> A5A5A5
$themeColors.Colors(8).RGB  = 49407      # accent4  CC9700 -> FFC000
$themeColors.Colors(9).RGB  = 12874308   # accent5  4EB3CF -> 4472C4
$themeColors.Colors(10).RGB = 4697456    # accent6  378DA6 -> 70AD47
$themeColors.Colors(11).RGB = 12673797   # hlink    6B9F25 -> 0563C1
$themeColors.Colors(12).RGB = 7491477    # folHlink B26B02 -> 954F72

$p.Save()
